# DIARIO DI BORDO - update dated 01-12-25
# 1) Bump the auto "datetimeFigureOut" date placeholders that were cached
#    from the previous save (27/11/2025 -> 01/12/2025) across the slide
#    master, every slide layout, and the notes master.
# 2) Update the wording of one of the discussion questions on slide 3.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Placeholders.Count; $i++) {
        $ph = $shapes.Placeholders.Item($i)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            $ph.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- Slide master: "11/27/2025" -> "12/1/2025" ---
Set-DatePlaceholderText $p.SlideMaster.Shapes "12/1/2025"

# --- Every slide layout: "11/27/2025" -> "12/1/2025" ---
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "12/1/2025"
}

# --- Notes master: "27/11/2025" -> "01/12/2025" ---
Set-DatePlaceholderText $p.NotesMaster.Shapes "01/12/2025"

# --- Slide 3: reword the verbali question ---
$slide3 = $p.Slides.Item(3)
$group4 = $slide3.Shapes.Item(3)
$textBox5 = $group4.GroupItems.Item(1)
$run = $textBox5.TextFrame.TextRange.Paragraphs(1).Runs(1)
$run.Text = "Se il ruolo incaricato della stesura dei verbali è il responsabile, chi è incaricato della loro successiva verifica e approvazione?"
